# 22/03/18 TK; All the boring stuff, config files, enums, etc..
$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-Para([string]$innerXml) {
    return "<w:p " + $ns + ">" + $innerXml + "</w:p>"
}

# Helper: replace the contents of an existing paragraph (by 1-based index) in place.
function Set-ParagraphXml($index, $innerXml) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertXML((New-Para $innerXml))
}

# Helper: insert a brand-new paragraph immediately after the paragraph at $index.
# Returns the index of the newly inserted paragraph.
function Add-ParagraphAfter($index, $innerXml) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertParagraphAfter()
    $newIndex = $index + 1
    $p2 = $d.Paragraphs.Item($newIndex)
    $p2.Range.InsertXML((New-Para $innerXml))
    return $newIndex
}

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that currently sits at the end of paragraph 3
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn the "Eclipse Generated "Comment Templates":" heading (paragraph 5)
#    into "Eclipse Perspectives:", then add two new paragraphs after it.
# ---------------------------------------------------------------------------
Set-ParagraphXml 5 '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Eclipse Perspectives:</w:t></w:r>'

$idx = Add-ParagraphAfter 5 '<w:r><w:t>Has everyone used / does everyone understand Eclipse perspectives?</w:t></w:r>'

$sourceInner = '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>'
$sourceInner += '<w:r><w:t xml:space="preserve">Eclipse </w:t></w:r>'
$sourceInner += '<w:r><w:t>– Sour</w:t></w:r>'
$sourceInner += '<w:r><w:t>c</w:t></w:r>'
$sourceInner += '<w:r><w:t xml:space="preserve">e - </w:t></w:r>'
$sourceInner += '<w:r><w:t xml:space="preserve">Generate </w:t></w:r>'
$sourceInner += '<w:r><w:t>Element Comments</w:t></w:r>'
$sourceInner += '<w:r><w:t>:</w:t></w:r>'
$idx = Add-ParagraphAfter $idx $sourceInner

# ---------------------------------------------------------------------------
# 3) Append the new "Refactor" / "Debugging" / "Resources Vs FileSystem"
#    sections right after the "General – Editors…" bullet point, i.e. right
#    before the trailing empty paragraph.
# ---------------------------------------------------------------------------
$lastBulletIndex = $d.Paragraphs.Count - 1
$idx = $lastBulletIndex

$idx = Add-ParagraphAfter $idx '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Refactor</w:t></w:r>'
$idx = Add-ParagraphAfter $idx '<w:r><w:t>Does everyone know how to refactor</w:t></w:r><w:r><w:t>?</w:t></w:r>'
$idx = Add-ParagraphAfter $idx '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Refactor – Extract Method</w:t></w:r>'

$extractInner = '<w:r><w:t>Do</w:t></w:r>'
$extractInner += '<w:r><w:t>es everyone know how to extract</w:t></w:r>'
$extractInner += '<w:r><w:t>?</w:t></w:r>'
$idx = Add-ParagraphAfter $idx $extractInner

$idx = Add-ParagraphAfter $idx '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Debugging:</w:t></w:r>'
$idx = Add-ParagraphAfter $idx '<w:r><w:t>Does everyone know how to run the debugger?</w:t></w:r>'

$resourcesHeadingInner = '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>'
$resourcesHeadingInner += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$resourcesHeadingInner += '<w:r><w:t xml:space="preserve">Resources Vs </w:t></w:r>'
$resourcesHeadingInner += '<w:proofErr w:type="spellStart"/>'
$resourcesHeadingInner += '<w:r><w:t>FileSystem</w:t></w:r>'
$resourcesHeadingInner += '<w:proofErr w:type="spellEnd"/>'
$resourcesHeadingInner += '<w:r><w:t>:</w:t></w:r>'
$idx = Add-ParagraphAfter $idx $resourcesHeadingInner

$resourcesBodyInner = '<w:r><w:t xml:space="preserve">Does everyone know </w:t></w:r>'
$resourcesBodyInner += '<w:r><w:t>the difference</w:t></w:r>'
$resourcesBodyInner += '<w:r><w:t>?</w:t></w:r>'
$resourcesBodyInner += '<w:r><w:t xml:space="preserve">  (Do I, really!? Lol…)</w:t></w:r>'
$idx = Add-ParagraphAfter $idx $resourcesBodyInner
